$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 222, shifting rows 222:270 down to 223:271.
$ws.Rows.Item(222).Insert()

# Populate the newly inserted row 222 with the new weekly record.
$ws.Cells.Item(222, 1).Value = 1
$ws.Cells.Item(222, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(222, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(222, 4).Value = 44641
$ws.Cells.Item(222, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(222, 5).Value = 15
$ws.Cells.Item(222, 6).Value = 100112043
$ws.Cells.Item(222, 7).Value = "Pepino ensalada"
$ws.Cells.Item(222, 8).Value = "Sin especificar"
$ws.Cells.Item(222, 9).Value = "Segunda"
$ws.Cells.Item(222, 10).Value = 120
$ws.Cells.Item(222, 11).Value = 12000
$ws.Cells.Item(222, 12).Value = 13000
$ws.Cells.Item(222, 13).Value = 12500
$ws.Cells.Item(222, 14).Value = "$/caja 100 unidades"
$ws.Cells.Item(222, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(222, 16).Value = 125
$ws.Cells.Item(222, 17).Value = 100
$ws.Cells.Item(222, 18).Value = "Hortaliza"
